$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 21: success-message test entry, styled like the other value cells (e.g. B9)
$ws.Range("A21").Value = "Success Message"
$ws.Range("B9").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B21").Value = "Success"

# Row 10: replace placeholder employee name and add new sample-data columns
$ws.Range("C10").Value = "a"
$ws.Range("D10").Value = "b"

# Row 5: "Add " -> "Add" (trailing space trimmed)
$ws.Range("B5").Value = "Add"

$ws.Range("B10").Value = "naveen  a"
$ws.Range("E10").Value = "s"

# Row 13: username value changed
$ws.Range("B13").Value = "aurasd312"

# Update the active selection to match the edited area
$ws.Range("B13").Select()
